# Rotate the data of rows 2, 3, 4 (columns A, B, E, F, G, H, Q, R) so that:
#   new row2 <- old row4
#   new row3 <- old row2
#   new row4 <- old row3
# All other columns are identical across these rows and remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values before overwriting anything.
# Use Value2 (not Value) since Value does not reliably return the
# underlying data through this COM-interop shim.
$orig = @{}
foreach ($r in 2..4) {
    $orig[$r] = @{
        A = $ws.Range("A$r").Value2
        B = $ws.Range("B$r").Value2
        E = $ws.Range("E$r").Value2
        F = $ws.Range("F$r").Value2
        G = $ws.Range("G$r").Value2
        H = $ws.Range("H$r").Value2
        Q = $ws.Range("Q$r").Value2
        R = $ws.Range("R$r").Value2
    }
}

# Mapping: destination row -> source row
$mapping = @{
    2 = 4
    3 = 2
    4 = 3
}

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    $data = $orig[$srcRow]

    $ws.Range("A$destRow").Value2 = $data.A
    $ws.Range("B$destRow").Value2 = $data.B
    $ws.Range("E$destRow").Value2 = $data.E
    $ws.Range("F$destRow").Value2 = $data.F
    $ws.Range("G$destRow").Value2 = $data.G
    $ws.Range("H$destRow").Value2 = $data.H
    $ws.Range("Q$destRow").Value2 = $data.Q
    $ws.Range("R$destRow").Value2 = $data.R
}
